$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 data (order chosen to match shared-string insertion order) ---
$ws.Range("B5").Value = "yes"
$ws.Range("D5").Value = "n.a."
$ws.Range("E5").Value = "global"
$ws.Range("F5").Value = "n.a."
$ws.Range("G5").Value = "general normative statement"
$ws.Range("H5").Value = "Calling for need to do something, no specific distribution. "
$ws.Range("C5").Value = "implementation, action"

# --- Re-apply the bold header font so styles.xml records family=2 for it ---
$ws.Range("A1:I1").Font.Name = "Calibri"

# --- Row height adjustments (from ~Windows 15pt default metrics to ~144/120 dpi re-layout) ---
$ws.Rows.Item(2).RowHeight = 144
$ws.Rows.Item(3).RowHeight = 115.2
$ws.Rows.Item(4).RowHeight = 57.6
$ws.Rows.Item(5).RowHeight = 72
$ws.Rows.Item(6).RowHeight = 129.6
$ws.Rows.Item(7).RowHeight = 57.6
$ws.Rows.Item(8).RowHeight = 172.8
$ws.Rows.Item(9).RowHeight = 187.2
$ws.Rows.Item(10).RowHeight = 142.2
$ws.Rows.Item(11).RowHeight = 115.2
$ws.Rows.Item(12).RowHeight = 259.2
$ws.Rows.Item(13).RowHeight = 86.4
$ws.Rows.Item(14).RowHeight = 100.8
$ws.Rows.Item(15).RowHeight = 43.2
$ws.Rows.Item(16).RowHeight = 43.2

# --- Update selection to D5 ---
$ws.Range("D5").Select()
